# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values for rows 2-21 change; everything else stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 5
    6  = 5
    7  = 4
    8  = 1
    9  = 6
    10 = 3
    11 = 4
    12 = 6
    13 = 3
    14 = 3
    15 = 1
    16 = 1
    17 = 1
    18 = 3
    19 = 2
    20 = 0
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
